$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (values unchanged, C1 stays "Links")
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Exclusions"
$ws.Range("C1").Value = "Links"

# New data for rows 2-13
$colA = @("Pat", "Em", "Bekah", "Kenny", "Jon", "Jo", "Chrissy", "Jessie", "Bill", "Mom", "Dad", "Andrew")
$colC = @("Em", "Pat", "Kenny", "Bekah", "Jo", "Jon", "Andrew", "Bill", "Jessie", "Dad", "Mom", "Chrissy")
$colB = @("Dad,Jessie", "Jessie,Kenny", "Mom,Em", "Jessie,Dad", "Emily,Bekah", "Bill,Mom", "Dad,Jo", "Jo,Jonathan,Chrissy", "Mom,Pat", "Chrissy,Bekah,Bill", "Kenny,Pat,Jon")

# Fill column A first (rows 2-13)
for ($i = 0; $i -lt $colA.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $colA[$i]
}

# Then fill column C (rows 2-13)
for ($i = 0; $i -lt $colC.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 3).Value = $colC[$i]
}

# Then fill column B (rows 2-12, row 13 has no exclusions)
for ($i = 0; $i -lt $colB.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $colB[$i]
}

$ws.Range("C18").Select()
